$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 633 entirely (the "木に登れば熊の危険..." post was removed),
# causing all subsequent rows to shift up by one.
$ws.Rows(633).Delete()
